# Auto-generated: apply crypto price/volume/hour updates + symbol-list shift (rows 49-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '302.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.15%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '18'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.19%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '18'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.105'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.00%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '18'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07703'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.50%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '18'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.420'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.47%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '18'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.615'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.49%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '18'

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '13.18%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '18'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1281'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '5.31%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '18'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1864'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.54%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '18'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09257'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.27%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '18'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04149'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.14%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '18'

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.61%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '18'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001277'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.18%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '18'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005760'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.54%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '18'

# Row 16
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '18'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.338'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.19%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '18'

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.89%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '18'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3340'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.73%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '18'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.082'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.55%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '18'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1370'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.66%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '18'

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '7.29%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '18'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04175'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.11%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '18'

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.58%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '18'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004399'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '13.53%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '18'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001349'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.76%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '18'

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '18'

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '18'

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '18'

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '18'

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '18'

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '18'

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '18'

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '18'

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '18'

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '18'

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '18'

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02521'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '4.67%'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '18'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05315'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.06%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '18'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.005926'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.09%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '18'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007716'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.94%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '18'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1351'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.03%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '18'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007348'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.13%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '18'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007530'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-7.03%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '18'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3030'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.02%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '18'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006681'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '7.01%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '18'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.12%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '18'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.04348'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-4.58%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '18'

# Row 49
$ws.Range("B49").Value = 'CryptobidCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.12%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '18'

# Row 50
$ws.Range("B50").Value = 'SpecialPowerGold'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.12%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '18'

# Row 51
$ws.Range("B51").Value = 'DigiFinexToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '--'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '--%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '18'
